# Generate Report for Handoff
# Updates the zh-cn and de-de report sheets: for the four ".md" rows whose
# handoff just completed (rows 4-7), bump Priority from "low" to "ht" and
# refresh the "Latest Handoff Datetime" timestamp.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $ws_zhcn.Cells.Item($r, 5).Value = "ht"
    $ws_zhcn.Cells.Item($r, 8).Value = "2016-08-29 12:33:27"
}

$ws_dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $ws_dede.Cells.Item($r, 5).Value = "ht"
    $ws_dede.Cells.Item($r, 8).Value = "2016-08-29 12:33:32"
}

# The Overview sheet mirrors the de-de "Latest HO Xliff Generate Date" value,
# so refresh it too to keep it in sync with the handoff run.
$ws_overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $ws_overview.Cells.Item($r, 7).Value = "2016-08-29 12:33:32"
}
